$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = '27.854.06'; E = '  +1.74%  ' }
    3 = @{ D = '1.810.32'; E = $null }
    4 = @{ D = '1.000'; E = '  -0.52%  ' }
    5 = @{ D = '336.95'; E = '  +0.36%  ' }
    6 = @{ D = '0.9962'; E = '  -0.62%  ' }
    7 = @{ D = '0.3922'; E = '  +3.68%  ' }
    8 = @{ D = '0.3485'; E = '  +1.71%  ' }
    9 = @{ D = '48.05'; E = '  -0.72%  ' }
    10 = @{ D = '1.199'; E = '  +0.11%  ' }
    11 = @{ D = '0.07571'; E = '  +1.46%  ' }
    12 = @{ D = '0.9982'; E = '  -0.43%  ' }
    13 = @{ D = '22.09'; E = '  +0.87%  ' }
    14 = @{ D = '6.511'; E = '  +0.82%  ' }
    15 = @{ D = '1.810.05'; E = '  +1.11%  ' }
    16 = @{ D = '7.189'; E = '  +2.36%  ' }
    17 = @{ D = '0.00001106'; E = '  +1.29%  ' }
    18 = @{ D = '0.06694'; E = '  +1.18%  ' }
    19 = @{ D = '85.18'; E = '  +0.99%  ' }
    20 = @{ D = '0.9950'; E = '  -0.69%  ' }
    21 = @{ D = '17.86'; E = '  +3.18%  ' }
    22 = @{ D = '6.564'; E = '  +1.71%  ' }
    23 = @{ D = '27.864.10'; E = '  +1.89%  ' }
    24 = @{ D = '12.85'; E = '  +2.65%  ' }
    25 = @{ D = '2.413'; E = '  -1.30%  ' }
    26 = @{ D = '2.549'; E = '  -0.32%  ' }
    27 = @{ D = '1.474'; E = '  +1.79%  ' }
    28 = @{ D = '21.28'; E = '  -0.36%  ' }
    29 = @{ D = '154.72'; E = '  +2.81%  ' }
    30 = @{ D = '2.017.73'; E = '  +1.27%  ' }
    31 = @{ D = '135.60'; E = '  +1.88%  ' }
    32 = @{ D = '4.037'; E = '  -0.39%  ' }
    33 = @{ D = '6.112'; E = '  +0.21%  ' }
    34 = @{ D = '0.08841'; E = '  +2.34%  ' }
    35 = @{ D = '13.26'; E = '  +0.16%  ' }
    36 = @{ D = '5.528'; E = '  +2.12%  ' }
    37 = @{ D = '0.02429'; E = '  +3.91%  ' }
    38 = @{ D = '0.06552'; E = '  +3.17%  ' }
    39 = @{ D = '0.6904'; E = '  +0.60%  ' }
    40 = @{ D = '1.612'; E = '  -3.39%  ' }
    41 = @{ D = '0.2224'; E = '  +1.37%  ' }
    42 = @{ D = '1.268'; E = '  +0.18%  ' }
    43 = @{ D = '8.574'; E = '  -2.55%  ' }
    44 = @{ D = '14.75'; E = '  +2.86%  ' }
    45 = @{ D = '0.6550'; E = '  +1.91%  ' }
    46 = @{ D = '0.9961'; E = '  -0.60%  ' }
    47 = @{ D = '3.865'; E = '  +0.43%  ' }
    48 = @{ D = '2.157'; E = '  +2.17%  ' }
    49 = @{ D = '132.42'; E = '  +2.48%  ' }
    50 = @{ D = '0.07201'; E = '  +0.06%  ' }
    51 = @{ D = '80.82'; E = '  +1.96%  ' }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item([int]$row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cellE = $ws.Cells.Item([int]$row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $u.E
    }
}
